$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Range("C1").Value = "Abreviatura"

# ---- Data: set numeric codes first (while cells are still General-formatted) ----
$ws.Range("A2").Value = 12
$ws.Range("A3").Value = 46
$ws.Range("A5").Value = 28

# ---- Text values for Nombre / Abreviatura columns ----
$ws.Range("B2").Value = "Castellón"
$ws.Range("C2").Value = "CAS"

$ws.Range("B3").Value = "Valencia"
$ws.Range("C3").Value = "VAL"

$ws.Range("B4").Value = "Barcelona"
$ws.Range("C4").Value = "BAR"

$ws.Range("B5").Value = "Madrid"
$ws.Range("C5").Value = "MAD"

# ---- A4 ("08") must be stored AS TEXT, so format it as Text before writing ----
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "08"

# ---- Apply the Text number format to the whole used range (A1:C5) ----
# (doesn't retroactively change already-stored numeric cells to text)
$ws.Range("A1:C5").NumberFormat = "@"

# ---- Selection shown in the saved file ----
$ws.Range("B9").Select() | Out-Null
